$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Sheet 1 - Crash Driver Report"

# Update the title cell (also updates the shared string "CrashDriverReport Complete")
$ws.Range("A1").Value = "Crash Driver Report"

# Row 20 was an exact duplicate of row 19 ("Crash" / "Crash Location Coordinates" / ...).
# Remove the duplicate row so the "Charge" rows that followed shift up one (old rows
# 21-24 become new rows 20-23).
$ws.Rows("20").Delete()
